$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.393.20"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.713.24"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5305"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06695"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2669"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -3.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07689"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.512"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").Value = "1.947.39"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "1.707.88"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5833"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "0.0₅8241"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "27.376.89"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "222.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.11%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.635"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.013"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.689"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1210"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.256"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05380"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.293"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.77%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.469"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.446"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.642"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.868"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9524"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.392"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5854"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01637"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.85%  "
$ws.Range("D40").Value = "1.094.67"
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.799"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8426"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").Value = "1.855.16"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").Value = "0.0₈117"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4537"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.006"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.131"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.11%  "
$ws.Range("E51").Value = "  -0.28%  "

Write-Output "done"
